{"js": "// 1. Update the title date from 27-06-23 to 11-07-23\nconst titleResults = context.document.body.search(\"WISC-IV rapport 27-06-23\", { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nfor (const r of titleResults.items) {\n  r.insertText(\"WISC-IV rapport 11-07-23\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. Change every occurrence of \"gennemsnitligt\" to \"nedre del af gennemsnittet\".\n//    This covers both the narrative paragraph text (\"... hvilket er gennemsnitligt.\")\n//    for the VSI and FHI indices, and the corresponding summary-table cells that\n//    contain just the word \"gennemsnitligt\" on its own.\nconst wordResults = context.document.body.search(\"gennemsnitligt\", { matchCase: true });\nwordResults.load(\"text\");\nawait context.sync();\nfor (const r of wordResults.items) {\n  r.insertText(\"nedre del af gennemsnittet\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the title date from 27-06-23 to 11-07-23\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"WISC-IV rapport 27-06-23\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue,\n    $false,\n    \"WISC-IV rapport 11-07-23\",\n    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceAll\n)\n\n# 2. Change every occurrence of \"gennemsnitligt\" to \"nedre del af gennemsnittet\".\n#    This covers both the narrative paragraph text (\"... hvilket er gennemsnitligt.\")\n#    for the VSI and FHI indices, and the corresponding summary-table cells that\n#    contain just the word \"gennemsnitligt\" on its own.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    \"gennemsnitligt\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue,\n    $false,\n    \"nedre del af gennemsnittet\",\n    [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceAll\n)\n\nWrite-Output \"done\"\n"}
